$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) column C for rows 2-14 moves forward by one day:
# from serial date 45174 (2023-09-05) to 45175 (2023-09-06).
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45175
}
